$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.660.84'
$ws.Range('E2').Value = '  +3.09%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.192.66'
$ws.Range('E3').Value = '  +0.74%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.32'
$ws.Range('E5').Value = '  +2.54%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '81.80'
$ws.Range('E6').Value = '  +11.62%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.624'
$ws.Range('E7').Value = '  +2.52%  '

$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.592'
$ws.Range('E9').Value = '  +1.97%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.27'
$ws.Range('E10').Value = '  +7.84%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0918'
$ws.Range('E11').Value = '  +0.64%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.96'
$ws.Range('E12').Value = '  +3.25%  '

$ws.Range('E13').Value = '  +2.14%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.516.70'
$ws.Range('E14').Value = '  +0.50%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.20'
$ws.Range('E15').Value = '  +0.51%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.193.28'
$ws.Range('E16').Value = '  +1.74%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.776'
$ws.Range('E17').Value = '  +1.02%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.560.86'
$ws.Range('E18').Value = '  +3.10%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000103'
$ws.Range('E19').Value = '  +0.93%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.02'
$ws.Range('E20').Value = '  -0.76%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.91'
$ws.Range('E21').Value = '  +0.97%  '

$ws.Range('E22').Value = '  +12.92%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.11'
$ws.Range('E23').Value = '  +1.32%  '

$ws.Range('E24').Value = '  -5.72%  '

$ws.Range('E25').Value = '  +0.19%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '42.07'
$ws.Range('E26').Value = '  +14.45%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.68'
$ws.Range('E27').Value = '  +2.46%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.37'
$ws.Range('E28').Value = '  -0.48%  '

$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('E29').Value = '  +2.89%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.19'
$ws.Range('E30').Value = '  -1.42%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.06'
$ws.Range('E31').Value = '  +1.34%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.41'
$ws.Range('E32').Value = '  +2.27%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0869'
$ws.Range('E33').Value = '  +7.61%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.29'
$ws.Range('E34').Value = '  +4.11%  '

$ws.Range('E35').Value = '  +6.38%  '

$ws.Range('E36').Value = '  +1.63%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.46'
$ws.Range('E37').Value = '  +5.85%  '

$ws.Range('E38').Value = '  +4.69%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.13'
$ws.Range('E39').Value = '  +12.11%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.85'
$ws.Range('E40').Value = '  +16.37%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.09'
$ws.Range('E41').Value = '  +2.51%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '62.66'
$ws.Range('E42').Value = '  +6.27%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.43'
$ws.Range('E43').Value = '  +6.40%  '

$ws.Range('E44').Value = '  +1.61%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.81'
$ws.Range('E45').Value = '  -0.52%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0980'
$ws.Range('E46').Value = '  +0.60%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.19'
$ws.Range('E47').Value = '  -0.07%  '

$ws.Range('E48').Value = '  +4.70%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.11'
$ws.Range('E49').Value = '  +2.44%  '

$ws.Range('E50').Value = '  +27.13%  '

$ws.Range('E51').Value = '  -5.87%  '
